$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-blank timelog rows 89-92 ------------------------
# These rows existed but were empty placeholders; the commit populates them
# with real "exposé"-style thesis-chapter work entries (date / from / to /
# duration-formula / activity), mirroring the layout of the rows above them.

# row 89: 2015-04-09, 16:00 - 19:30
$ws.Range("A89").Value = 42103
$ws.Range("B89").Value = 0.66666666666666663
$ws.Range("C89").Value = 0.8125
$ws.Range("D89").Formula = "=C89-B89"
$ws.Range("E89").Value = "thesis chapter introduction"

# row 90: 2015-04-10, 10:00 - 13:00
$ws.Range("A90").Value = 42104
$ws.Range("B90").Value = 0.41666666666666669
$ws.Range("C90").Value = 0.54166666666666663
$ws.Range("D90").Formula = "=C90-B90"
$ws.Range("E90").Value = "thesis chapter introduction"

# row 91: 2015-04-15, 11:00 - 14:00
$ws.Range("A91").Value = 42109
$ws.Range("B91").Value = 0.45833333333333331
$ws.Range("C91").Value = 0.58333333333333337
$ws.Range("D91").Formula = "=C91-B91"
$ws.Range("E91").Value = "thesis chapter introduction"

# row 92: 2015-04-15, 16:30 - 18:30
$ws.Range("A92").Value = 42109
$ws.Range("B92").Value = 0.6875
$ws.Range("C92").Value = 0.77083333333333337
$ws.Range("D92").Formula = "=C92-B92"
$ws.Range("E92").Value = "thesis chapter introduction"

# Row 93 / the D94 total stay formula-driven (=SUM(D2:D93)) and pick up the
# new rows automatically on recalc.

# --- Cursor/selection moves from E91 to I88 --------------------------------
$ws.Range("I88").Select() | Out-Null
